$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for column AB (one new day of COVID case counts), continuing
# directly on from column AA.
$newValues = @(267, 130, 141, 200, 226, 194, 931, 99, 478, 534, 4, 41, 457, 4)

# Copy the formatting of column AA into AB so the new column matches the
# existing look (borders, number format, etc.), then fill in the new values.
$ws.Range("AA1:AA14").Copy() | Out-Null
$ws.Range("AB1:AB14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 28).Value = $newValues[$i]
}

# Update the view so the new column is visible/selected, matching the
# author's last on-screen state (scrolled right, whole table selected).
$ws.Application.CutCopyMode = $false
$activeWindow = $ws.Application.ActiveWindow
if ($activeWindow) {
    $activeWindow.ScrollColumn = 24
}
$ws.Range("A1:AB14").Select() | Out-Null
